$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 320. Excel shifts the existing rows 320:421
# down to 322:423 (preserving all their data/formatting), which already
# reproduces the entire downstream diff (each row's content simply moves
# down by two rows, with the former last pair ending up at rows 422:423).
$ws.Rows("320:321").Insert()

# Populate the two new rows (320 = "Primera", 321 = "Segunda") with this
# week's fresh data point. Static columns (A,B,C,E,F,G,H,I,N,O,Q,R) are the
# same as every other row in this block, so copy them from the row that
# used to be first (now row 322/323) and only change the columns that the
# diff actually shows changing (D, J, K, L, M, P).
$ws.Range("A320:R320").Value = $ws.Range("A322:R322").Value2
$ws.Range("A321:R321").Value = $ws.Range("A323:R323").Value2

$ws.Range("D320").Value = 44722
$ws.Range("J320").Value = 2500
$ws.Range("K320").Value = 600
$ws.Range("L320").Value = 700
$ws.Range("M320").Value = 650
$ws.Range("P320").Value = 325

$ws.Range("D321").Value = 44722
$ws.Range("J321").Value = 1500
$ws.Range("K321").Value = 500
$ws.Range("L321").Value = 550
$ws.Range("M321").Value = 525
$ws.Range("P321").Value = 262
